$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF column (F) values with repulled data
$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = -3
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = -2

# Row 23 also has an updated PC (D) value
$ws.Range("D23").Value = 97
$ws.Range("F23").Value = -2

$ws.Range("F24").Value = -4
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = -1
$ws.Range("F27").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("F30").Value = -2
